$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" column (S) is being appended after the existing "2021"
# column (R). Copy R2:R6's formatting into S2:S6 so the new column
# matches the look of the existing year columns, then fill in the
# 2022 figures.
$ws.Range("R2:R6").Copy()
[void]$ws.Range("S2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# The author's saved cursor position after the edit.
[void]$ws.Range("C19").Select()
